# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# These same events appear on both the "展览" sheet and the aggregated
# "全部类型" sheet, so both need to be updated in lockstep.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1162
$ws1.Range("F4").Value  = 16591
$ws1.Range("F11").Value = 11545
$ws1.Range("F13").Value = 1225
$ws1.Range("F14").Value = 4570
$ws1.Range("F15").Value = 402
$ws1.Range("F18").Value = 871

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1162
$ws4.Range("F5").Value  = 16591
$ws4.Range("F14").Value = 11545
$ws4.Range("F16").Value = 1225
$ws4.Range("F17").Value = 4570
$ws4.Range("F18").Value = 402
$ws4.Range("F21").Value = 871
